$wb = $excel.ActiveWorkbook

# Rename the only worksheet from "Tabelle1" to "Sheet1" (matches pandas' default
# sheet naming convention referenced in the commit message).
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"
